$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These values must remain stored as text (matching the original
# inline-string cells), not be reinterpreted as numbers. A leading
# apostrophe forces text entry; resetting the style back to "Normal"
# afterwards avoids leaving a stray quote-prefix style applied to the
# cell (keeping formatting identical to the original).
$ws.Range("B9").Value = "'45.74"
$ws.Range("B9").Style = "Normal"

$ws.Range("B10").Value = "'48.79"
$ws.Range("B10").Style = "Normal"

$ws.Range("B11").Value = "'0.32"
$ws.Range("B11").Style = "Normal"

$ws.Range("B12").Value = "'93.02"
$ws.Range("B12").Style = "Normal"

$ws.Range("B13").Value = "'6.94"
$ws.Range("B13").Style = "Normal"
